# Assignment1.xlsx edit script
# - Updates several data points on the "PO List" sheet (these feed the
#   dynamic-array FILTER formula on the "Assignment" sheet).
# - Shrinks the FILTER array-formula spill range on "Assignment"!C13 from
#   C13:C16 down to C13:C15 (one fewer name now satisfies the filter).
# - Makes "PO List" the active/selected sheet instead of "PO GBW List".

$wb  = $excel.ActiveWorkbook
$wsAssignment = $wb.Worksheets.Item("Assignment")
$wsPoList     = $wb.Worksheets.Item("PO List")

# ---------------------------------------------------------------------
# 1. Update the underlying data on the "PO List" sheet.
# ---------------------------------------------------------------------
$wsPoList.Range("Y4").Value  = 6

$wsPoList.Range("N7").Value  = 3

$wsPoList.Range("N9").Value  = 8
$wsPoList.Range("Y9").Value  = 4

$wsPoList.Range("K11").Value = 8
$wsPoList.Range("L11").Value = 8
$wsPoList.Range("M11").Value = 44951
$wsPoList.Range("N11").Value = 1
$wsPoList.Range("S11").Value = 3
$wsPoList.Range("T11").Value = 3
$wsPoList.Range("U11").Value = 44951
$wsPoList.Range("V11").Value = 9
$wsPoList.Range("W11").Value = 9
$wsPoList.Range("X11").Value = 44951
$wsPoList.Range("Y11").Value = 1

$wsPoList.Range("N12").Value = 10
$wsPoList.Range("Y12").Value = 5

$wsPoList.Range("N13").Value = 9
$wsPoList.Range("O13").Value = 1
$wsPoList.Range("P13").Value = 1
$wsPoList.Range("Q13").Value = 44946
$wsPoList.Range("R13").Value = 1
$wsPoList.Range("S13").Value = 1
$wsPoList.Range("T13").Value = 1
$wsPoList.Range("U13").Value = 44946

$wsPoList.Range("R14").Value = 25

$wsPoList.Range("Y16").Value = 7

$wsPoList.Range("O17").Value = 2
$wsPoList.Range("P17").Value = 2
$wsPoList.Range("Q17").Value = 44120
$wsPoList.Range("R17").Value = 23
$wsPoList.Range("S17").Value = 3
$wsPoList.Range("T17").Value = 3
$wsPoList.Range("U17").Value = 44809

$wsPoList.Range("N19").Value = 6
$wsPoList.Range("R19").Value = 25
$wsPoList.Range("Y19").Value = 3

$wsPoList.Range("N20").Value = 3

$wsPoList.Range("R22").Value = 24

$wsPoList.Range("N24").Value = 11

$wsPoList.Range("N25").Value = 7

$wsPoList.Range("N27").Value = 5

$wsPoList.Range("N28").Value = 2
$wsPoList.Range("R28").Value = 25
$wsPoList.Range("Y28").Value = 2

# ---------------------------------------------------------------------
# 2. Re-enter the FILTER array formula on "Assignment" so its spill
#    range shrinks from C13:C16 to C13:C15 (now that the data above
#    leaves one fewer matching row).
# ---------------------------------------------------------------------
$wsAssignment.Range("C13:C16").ClearContents()
$wsAssignment.Range("C13:C15").FormulaArray = " _xlfn._xlws.FILTER('PO List'!`$B`$3:`$B`$29, 'PO List'!`$O`$3:`$O`$29 = MAX(_xlfn._xlws.FILTER('PO List'!`$O`$3:`$O`$29,  ('PO List'!`$W`$3:`$W`$29=1) * ('PO List'!F3:F29>=3))) * ('PO List'!`$W`$3:`$W`$29 = 1) * ('PO List'!F3:F29>=3))"

# ---------------------------------------------------------------------
# 3. Switch the active/selected sheet from "PO GBW List" to "PO List".
# ---------------------------------------------------------------------
$wsPoList.Activate()
